$wb = $excel.ActiveWorkbook

# Update the "Functions" sheet: Battleship, Player and Randomizer functions
# (rows 34-58) are untestable without mocking, so their Tester/Status columns
# move from "?"/owner-name (and a blank Status) to "N/A".
$functions = $wb.Worksheets.Item("Functions")
for ($r = 34; $r -le 58; $r++) {
    $functions.Cells.Item($r, 5).Value = "N/A"
    $functions.Cells.Item($r, 6).Value = "N/A"
}

# Update the "SW Units" sheet: Player is no longer testable, and needs mocking.
$swUnits = $wb.Worksheets.Item("SW Units")
$swUnits.Range("B5").Value = "No"
$swUnits.Range("C5").Value = "Mocking required."
